$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("N3").Value = 1.36
$ws.Range("P3").Value = 1.36

# Row 6 updates (F6:AO6)
$ws.Range("F6").Value = 4.7
$ws.Range("G6").Value = 7.6
$ws.Range("H6").Value = 1.57
$ws.Range("I6").Value = 1.74
$ws.Range("J6").Value = 3.7
$ws.Range("K6").Value = 5.6
$ws.Range("L6").Value = 1.01
$ws.Range("M6").Value = 1.04
$ws.Range("N6").Value = 3.9
$ws.Range("O6").Value = 1.23
$ws.Range("P6").Value = 2.12
$ws.Range("Q6").Value = 1.68
$ws.Range("R6").Value = 1.44
$ws.Range("S6").Value = 2.66
$ws.Range("T6").Value = 1.74
$ws.Range("U6").Value = 2.02
$ws.Range("V6").Value = 2.34
$ws.Range("W6").Value = 1.17
$ws.Range("X6").Value = 990
$ws.Range("Y6").Value = 990
$ws.Range("Z6").Value = 980
$ws.Range("AA6").Value = 980
$ws.Range("AB6").Value = 990
$ws.Range("AC6").Value = 990
$ws.Range("AD6").Value = 990
$ws.Range("AE6").Value = 980
$ws.Range("AF6").Value = 60
$ws.Range("AG6").Value = 990
$ws.Range("AH6").Value = 990
$ws.Range("AI6").Value = 980
$ws.Range("AJ6").Value = 1000
$ws.Range("AK6").Value = 1000
$ws.Range("AL6").Value = 85
$ws.Range("AM6").Value = 1000
$ws.Range("AN6").Value = 100
$ws.Range("AO6").Value = 980

# Row 7 updates (F7:AO7)
$ws.Range("F7").Value = 3.05
$ws.Range("G7").Value = 3.6
$ws.Range("H7").Value = 2.04
$ws.Range("I7").Value = 2.28
$ws.Range("J7").Value = 3.95
$ws.Range("K7").Value = 4.7
$ws.Range("L7").Value = 1.01
$ws.Range("M7").Value = 1.03
$ws.Range("N7").Value = 5.8
$ws.Range("O7").Value = 1.14
$ws.Range("P7").Value = 2.62
$ws.Range("Q7").Value = 1.49
$ws.Range("R7").Value = 1.66
$ws.Range("S7").Value = 2.22
$ws.Range("T7").Value = 1.49
$ws.Range("U7").Value = 2.62
$ws.Range("V7").Value = 1.79
$ws.Range("W7").Value = 1.39
$ws.Range("X7").Value = 36
$ws.Range("Y7").Value = 19.5
$ws.Range("Z7").Value = 22
$ws.Range("AA7").Value = 34
$ws.Range("AB7").Value = 22
$ws.Range("AC7").Value = 13
$ws.Range("AD7").Value = 13.5
$ws.Range("AE7").Value = 21
$ws.Range("AF7").Value = 36
$ws.Range("AG7").Value = 18.5
$ws.Range("AH7").Value = 18
$ws.Range("AI7").Value = 32
$ws.Range("AJ7").Value = 65
$ws.Range("AK7").Value = 38
$ws.Range("AL7").Value = 38
$ws.Range("AM7").Value = 65
$ws.Range("AN7").Value = 24
$ws.Range("AO7").Value = 11.5

# Row 8 update
$ws.Range("H8").Value = 3.05
